$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CNN 2D")
$ws.Activate()

$data = @(
    @(3, 88.130527734756399, 0.27344635745206602, 67.290010929107595),
    @(4, 88.920247554778996, 0.25774980586164398, 66.396825075149494),
    @(5, 83.821469545364295, 0.297402189457594, 65.558467149734497),
    @(6, 87.620651721954303, 0.28040545287405499, 65.826833248138399),
    @(7, 87.561362981796194, 0.29203820308600498, 65.761678695678697),
    @(8, 83.892619609832707, 0.288058755348577, 65.0974698066711),
    @(9, 88.011950254440293, 0.26516301425847799, 61.287345170974703),
    @(10, 87.803256511688204, 0.27637953875471399, 61.200113534927297),
    @(11, 87.563735246658297, 0.28503508429721502, 66.596321582794104),
    @(12, 87.203264236450195, 0.28297697417776302, 64.106638193130493),
    @(13, 87.947922945022498, 0.26649139156371499, 64.160237550735403),
    @(14, 88.524198532104407, 0.26339664288289399, 61.358152866363497),
    @(15, 83.712381124496403, 0.28441361926408898, 65.032004833221393),
    @(16, 83.736097812652503, 0.28139868792735001, 65.463674068450899),
    @(17, 87.848317623138399, 0.26422198295238403, 66.049296617507906),
    @(18, 83.264166116714406, 0.30746913967835998, 61.688814640045102),
    @(19, 86.653071641921997, 0.29631816906081898, 65.205819129943805),
    @(20, 88.711553812026906, 0.26160770693638902, 66.543475866317706),
    @(21, 87.089431285858097, 0.29469060519281698, 67.582562446594196),
    @(22, 82.941639423370304, 2.6234215325207999, 67.450929641723604),
    @(23, 87.440413236617999, 0.26581365415767599, 67.692631244659395),
    @(24, 87.053859233856201, 0.327350065479457, 67.456271886825505),
    @(25, 88.417482376098604, 0.249481786307755, 64.375463724136296),
    @(26, 87.826973199844304, 0.271586461314759, 60.734239101409898),
    @(27, 86.980342864990206, 0.29200259537853901, 66.988158702850299),
    @(28, 87.817490100860596, 0.25658139156344401, 61.1500530242919),
    @(29, 88.464915752410803, 0.25718509018241598, 60.492553949356001),
    @(30, 87.881517410278306, 0.26282472246713801, 66.510779142379704),
    @(31, 87.762939929962101, 0.27445424753006098, 66.378743171691895),
    @(32, 86.747932434082003, 0.29183078213038699, 61.351038455963099),
    @(33, 87.698912620544405, 0.273095712291332, 66.023054361343299),
    @(34, 88.0166947841644, 0.25742864137367399, 67.002029418945298),
    @(35, 88.078355789184499, 0.25665540870306502, 61.413606643676701),
    @(36, 87.722629308700505, 0.25741935373843799, 67.168449163436804),
    @(37, 83.209615945815997, 0.328430703784232, 66.8585751056671),
    @(38, 87.698912620544405, 0.27151091939571698, 67.386412620544405),
    @(39, 87.241208553314195, 0.26767146111504297, 63.078752994537297),
    @(40, 87.696540355682302, 0.29575851041206302, 66.118104219436603),
    @(41, 88.147127628326402, 0.268626176879037, 60.746504306793199),
    @(42, 83.366137742996202, 0.31244047473066799, 66.8369717597961),
    @(43, 84.260201454162598, 0.28698463674673402, 66.3361496925354),
    @(44, 87.623023986816406, 0.26719601728273101, 61.371240615844698),
    @(45, 88.657009601593003, 0.25705344326551899, 64.951596975326495),
    @(46, 87.577962875366197, 0.276600441735673, 60.868054389953599),
    @(47, 87.224608659744206, 0.264391790864409, 66.776542425155597),
    @(48, 86.629354953765798, 0.30728197478090602, 65.724978685379),
    @(49, 87.684684991836505, 0.28741314188365003, 66.127233982086096),
    @(50, 87.649106979370103, 0.30061774726076301, 64.955434799194293),
    @(51, 88.358193635940495, 0.26374111377715398, 61.443469762802103)
)

foreach ($item in $data) {
    $rowNum = $item[0]
    $ws.Cells.Item($rowNum, 25).Value = $item[1]   # Y
    $ws.Cells.Item($rowNum, 26).Value = $item[2]   # Z
    $ws.Cells.Item($rowNum, 27).Value = $item[3]   # AA
}

$ws.Range("U11").Select()
